$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (id, number, name, idno, phone, email, adr, password, created_at).
# All lecturer data rows shift up by one (EZANA001 now on row 1, ... EZANA028 now on row 28).
$ws.Rows.Item(1).Delete()

# Give every lecturer row (from EZANA003 / new row 3 onward) its own unique e-mail address
# instead of the previously shared "lec002@ezana.org" placeholder. Row 10 keeps the source
# data's "lec0010@ezana.org" typo.
$emails = @{}
$emails[3]  = "lec003@ezana.org"
$emails[4]  = "lec004@ezana.org"
$emails[5]  = "lec005@ezana.org"
$emails[6]  = "lec006@ezana.org"
$emails[7]  = "lec007@ezana.org"
$emails[8]  = "lec008@ezana.org"
$emails[9]  = "lec009@ezana.org"
$emails[10] = "lec0010@ezana.org"
$emails[11] = "lec011@ezana.org"
$emails[12] = "lec012@ezana.org"
$emails[13] = "lec013@ezana.org"
$emails[14] = "lec014@ezana.org"
$emails[15] = "lec015@ezana.org"
$emails[16] = "lec016@ezana.org"
$emails[17] = "lec017@ezana.org"
$emails[18] = "lec018@ezana.org"
$emails[19] = "lec019@ezana.org"
$emails[20] = "lec020@ezana.org"
$emails[21] = "lec021@ezana.org"
$emails[22] = "lec022@ezana.org"
$emails[23] = "lec023@ezana.org"
$emails[24] = "lec024@ezana.org"
$emails[25] = "lec025@ezana.org"
$emails[26] = "lec026@ezana.org"
$emails[27] = "lec027@ezana.org"
$emails[28] = "lec028@ezana.org"

foreach ($row in $emails.Keys) {
    $ws.Cells.Item($row, 6).Value = $emails[$row]
}

# Match the author's final selection.
$ws.Range("F29").Select()
